# Update the cryptocurrency price/volume table (columns D and E) with refreshed
# figures pulled by the scheduled GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.606.04"
$ws.Range("E2").Value = "  +1.15%  "
$ws.Range("D3").Value = "3.394.57"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.36"
$ws.Range("E5").Value = "  +0.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.44"
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.73"
$ws.Range("E9").Value = "  +2.63%  "
$ws.Range("E10").Value = "  -0.49%  "
$ws.Range("E11").Value = "  -1.61%  "
$ws.Range("D12").Value = "3.975.27"
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.47"
$ws.Range("E13").Value = "  +0.99%  "
$ws.Range("E14").Value = "  +0.15%  "
$ws.Range("D15").Value = "3.407.40"
$ws.Range("E15").Value = "  +0.55%  "
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("D17").Value = "61.636.27"
$ws.Range("E17").Value = "  +1.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.17"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.70"
$ws.Range("E19").Value = "  -0.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.98"
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "391.26"
$ws.Range("E21").Value = "  +1.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "75.49"
$ws.Range("E22").Value = "  +1.40%  "
$ws.Range("E23").Value = "  -0.55%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("E25").Value = "  -3.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.194"
$ws.Range("E26").Value = "  +8.32%  "
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.27"
$ws.Range("E28").Value = "  -1.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.03"
$ws.Range("E29").Value = "  +0.74%  "
$ws.Range("E30").Value = "  +0.33%  "
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("E32").Value = "  -3.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.37"
$ws.Range("E33").Value = "  -0.42%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.94"
$ws.Range("E34").Value = "  -0.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "168.16"
$ws.Range("E35").Value = "  +0.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.07"
$ws.Range("E36").Value = "  +1.85%  "
$ws.Range("D37").Value = "3.430.14"
$ws.Range("E37").Value = "  +0.41%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.47"
$ws.Range("E38").Value = "  -0.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0772"
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "25.95"
$ws.Range("E40").Value = "  -6.37%  "
$ws.Range("E41").Value = "  +0.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.43"
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.66"
$ws.Range("E43").Value = "  -0.41%  "
$ws.Range("E44").Value = "  +1.59%  "
$ws.Range("D45").Value = "2.472.89"
$ws.Range("E45").Value = "  -0.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "23.07"
$ws.Range("E46").Value = "  +0.32%  "
$ws.Range("E47").Value = "  -2.07%  "
$ws.Range("E48").Value = "  +0.06%  "
$ws.Range("E49").Value = "  -1.60%  "
$ws.Range("E50").Value = "  -0.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.206"
$ws.Range("E51").Value = "  -1.45%  "
